$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the 6 stray empty paragraphs (self-closing <w:p/>) that sit
#    between the heading/edition/body paragraphs near the top of the
#    document.  They are paragraphs 6, 8, 10, 12, 14 and 16 (1-based) in the
#    original document -- delete from the bottom up so indices stay valid.
# ---------------------------------------------------------------------------
$emptyIndices = @(16, 14, 12, 10, 8, 6)
foreach ($idx in $emptyIndices) {
    $p = $d.Paragraphs.Item($idx)
    if ($p.Range.Text.Trim().Length -eq 0) {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Style tweaks (styles.xml)
# ---------------------------------------------------------------------------

# Normal: after=240  ->  before=240, after=240
$sNormal = $d.Styles.Item("Normal")
$sNormal.ParagraphFormat.SpaceBefore = 12
$sNormal.ParagraphFormat.SpaceAfter = 12

# Heading 3: before=240 after=0  ->  before=360 (after now inherited = 240)
$sHeading3 = $d.Styles.Item("Heading 3")
$sHeading3.ParagraphFormat.SpaceBefore = 18
$sHeading3.ParagraphFormat.SpaceAfter = 12

# Heading 6: before=240 after=60  ->  after=60 (before now inherited = 240)
$sHeading6 = $d.Styles.Item("Heading 6")
$sHeading6.ParagraphFormat.SpaceBefore = 12
$sHeading6.ParagraphFormat.SpaceAfter = 3

# List 1: before=120  ->  before=240 after=240 ; contextualSpacing removed
$sList1 = $d.Styles.Item("List 1")
$sList1.ParagraphFormat.SpaceBefore = 12
$sList1.ParagraphFormat.SpaceAfter = 12
$sList1.NoSpaceBetweenParagraphsOfSameStyle = $false

# List 2/3/4: spacing + contextualSpacing removed -> inherits Normal (240/240)
$sList2 = $d.Styles.Item("List 2")
$sList2.ParagraphFormat.SpaceBefore = 12
$sList2.ParagraphFormat.SpaceAfter = 12
$sList2.NoSpaceBetweenParagraphsOfSameStyle = $false

$sList3 = $d.Styles.Item("List 3")
$sList3.ParagraphFormat.SpaceBefore = 12
$sList3.ParagraphFormat.SpaceAfter = 12
$sList3.NoSpaceBetweenParagraphsOfSameStyle = $false

$sList4 = $d.Styles.Item("List 4")
$sList4.ParagraphFormat.SpaceBefore = 12
$sList4.ParagraphFormat.SpaceAfter = 12
$sList4.NoSpaceBetweenParagraphsOfSameStyle = $false

# List 1_change: spacing + contextualSpacing removed -> inherits Normal (240/240)
$sList1change = $d.Styles.Item("List 1_change")
$sList1change.ParagraphFormat.SpaceBefore = 12
$sList1change.ParagraphFormat.SpaceAfter = 12
$sList1change.NoSpaceBetweenParagraphsOfSameStyle = $false

# List 1_change Char: drop the theme-color reference, keep plain black RGB
$sList1changeChar = $d.Styles.Item("List 1_change Char")
$sList1changeChar.Font.TextColor.RGB = 0

# TOC Heading: before=240 after=0 line=259/auto -> after=0 line=259/auto
# (before now inherited via Heading 1 -> Normal = 240)
$sTOCHeading = $d.Styles.Item("TOC Heading")
$sTOCHeading.ParagraphFormat.SpaceBefore = 12
$sTOCHeading.ParagraphFormat.SpaceAfter = 0
